$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 7: "TextBox 47" -> "[command commits address book]" becomes
# "[command commits order book]" split into the runs
# "[", "command commits ", "order book", "]"
$shpCommand = $s.Shapes.Item(7)
$origHeightCommand = $shpCommand.Height
$trCommand = $shpCommand.TextFrame.TextRange
$subCommand = $trCommand.Characters(18, 12)
$subCommand.Text = "order book"
# Restore the autofit textbox height PowerPoint nudges on text edit, so the
# shape geometry stays as close as possible to its original size.
$shpCommand.Height = $origHeightCommand

# Shape 8: "Rectangle: Rounded Corners 50" ->
# "Purge redundant states and then save address book to addressBookStateList "
# becomes
# "Purge redundant states and then save order book to orderBookStateList "
$shpPurge = $s.Shapes.Item(8)
$trPurge = $shpPurge.TextFrame.TextRange
$subPurge1 = $trPurge.Characters(38, 13)
$subPurge1.Text = "order book "
$subPurge2 = $trPurge.Characters(52, 20)
$subPurge2.Text = "orderBookStateList"
